$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Delete rows 15, 14, 5 (from bottom to top so row indices of rows
# still to be removed remain valid as each delete shifts subsequent rows up).
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(5).Delete()
